# Check_list.xlsx - add a new "Страница тематических цитат" (Thematic quotes
# page) section at the bottom of the checklist (rows 39-41), matching the
# pattern used by the other section headers already present on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122  # xlPasteFormats

# --- Row 39: new section header (merged A39:I39, same look as A2/A10/A16/A24) ---
$ws.Range("A2:I2").Copy()
$ws.Range("A39:I39").PasteSpecial($xlPasteFormats)
$ws.Range("A39").Value = "Страница тематических цитат"
$ws.Range("A39:I39").Merge()

# --- Row 40: first check under the new section (same style as A11/B11) ---
$ws.Range("A11").Copy()
$ws.Range("A40").PasteSpecial($xlPasteFormats)
$ws.Range("B11").Copy()
$ws.Range("B40").PasteSpecial($xlPasteFormats)
$ws.Range("A40").Value = "Переход на страницу тематических цитат"
$ws.Range("B40").Value = "Pass"

# --- Row 41: second check under the new section ---
$ws.Range("A11").Copy()
$ws.Range("A41").PasteSpecial($xlPasteFormats)
$ws.Range("B11").Copy()
$ws.Range("B41").PasteSpecial($xlPasteFormats)
$ws.Range("A41").Value = "Развернуть тематическую цитату"
$ws.Range("B41").Value = "Pass"

# --- Update the view's selection to match the edited area ---
$ws.Range("D43").Select()

$excel.CutCopyMode = $false
